# Update gh-pages to output generated at 456a3b4
#
# Refreshes the scraped "想去人数" (want-to-go) counts on several existing
# rows and appends a new event row (id 17, "合肥·运动番only-群青日和") to
# both the "展览" and "全部类型" sheets.

$xlPasteFormats = [int][Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # --- Refreshed "想去人数" counts for existing rows -------------------
    $ws.Range("F3").Value = 7241
    if ($name -eq "展览") {
        $ws.Range("F4").Value = 5393
    } else {
        $ws.Range("F4").Value = 5394
    }
    $ws.Range("F5").Value = 78
    $ws.Range("F10").Value = 81
    $ws.Range("F11").Value = 100
    $ws.Range("F12").Value = 200
    $ws.Range("F13").Value = 13
    $ws.Range("F14").Value = 644
    $ws.Range("F15").Value = 247

    # --- New row 18: 合肥·运动番only-群青日和 -----------------------------
    $ws.Cells.Item(18, 1).Value = 17
    $ws.Cells.Item(18, 2).Value = "'2024-06-01"
    $ws.Cells.Item(18, 3).Value = "合肥·运动番only-群青日和"
    $ws.Cells.Item(18, 4).Value = "金寨路287号 合肥明星运动公园"
    $ws.Cells.Item(18, 5).Value = "2024.06.01 09:30-06.01 17:30"
    $ws.Cells.Item(18, 6).Value = 0
    $ws.Cells.Item(18, 7).Value = 70
    $ws.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83058"
    $ws.Cells.Item(18, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/RcH1W6cK1710422301382.jpeg"

    # Column A uses the bold/bordered/centered style applied to the rest of
    # column A (same as row 17) - copy formats only, the value is already set.
    $ws.Range("A17").Copy() | Out-Null
    $ws.Range("A18").PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = 0
}
